$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '39.527.33'
$ws.Range("E2").Value = '  +1.83%  '

# Row 3
$ws.Range("D3").Value = '2.157.45'
$ws.Range("E3").Value = '  +3.08%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").Value = "'229.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.54%  '

# Row 6
$ws.Range("D6").Value = "'0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.33%  '

# Row 7
$ws.Range("E7").Value = '  +4.26%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("E9").Value = '  +2.64%  '

# Row 10
$ws.Range("D10").Value = "'0.0860"
$ws.Range("D10").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'16.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.98%  '

# Row 13
$ws.Range("D13").Value = '2.477.30'
$ws.Range("E13").Value = '  +3.13%  '

# Row 14
$ws.Range("D14").Value = "'22.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.86%  '

# Row 15
$ws.Range("D15").Value = "'0.819"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.73%  '

# Row 16
$ws.Range("D16").Value = "'5.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.00%  '

# Row 17
$ws.Range("D17").Value = '2.153.63'
$ws.Range("E17").Value = '  +2.94%  '

# Row 18
$ws.Range("D18").Value = '39.503.21'
$ws.Range("E18").Value = '  +2.08%  '

# Row 19
$ws.Range("D19").Value = "'72.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.14%  '

# Row 20
$ws.Range("D20").Value = "'6.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.91%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0855'
$ws.Range("E21").Value = '  +1.99%  '

# Row 22
$ws.Range("D22").Value = "'229.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.82%  '

# Row 23
$ws.Range("E23").Value = '  +0.01%  '

# Row 24
$ws.Range("E24").Value = '  -0.64%  '

# Row 25
$ws.Range("D25").Value = "'2.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.38%  '

# Row 26
$ws.Range("D26").Value = "'9.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.24%  '

# Row 27
$ws.Range("D27").Value = "'173.10"
$ws.Range("D27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'0.139"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.55%  '

# Row 29
$ws.Range("E29").Value = '  -1.70%  '

# Row 30
$ws.Range("E30").Value = '  +2.52%  '

# Row 31
$ws.Range("E31").Value = '  +8.72%  '

# Row 32
$ws.Range("E32").Value = '  +1.31%  '

# Row 33
$ws.Range("E33").Value = '  +3.63%  '

# Row 34
$ws.Range("E34").Value = '  +2.99%  '

# Row 35
$ws.Range("D35").Value = "'7.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.63%  '

# Row 36
$ws.Range("D36").Value = "'0.0624"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.79%  '

# Row 37
$ws.Range("E37").Value = '  +1.51%  '

# Row 38
$ws.Range("E38").Value = '  -0.04%  '

# Row 39
$ws.Range("E39").Value = '  +0.20%  '

# Row 40
$ws.Range("E40").Value = '  +0.82%  '

# Row 41
$ws.Range("E41").Value = '  +3.48%  '

# Row 42
$ws.Range("D42").Value = "'103.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.54%  '

# Row 43
$ws.Range("D43").Value = '1.537.37'
$ws.Range("E43").Value = '  -0.30%  '

# Row 44
$ws.Range("E44").Value = '  +6.27%  '

# Row 45
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = "'0.0928"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.67%  '

# Row 46
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = "'1.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.05%  '

# Row 47
$ws.Range("E47").Value = '  -0.34%  '

# Row 48
$ws.Range("D48").Value = "'7.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.94%  '

# Row 49
$ws.Range("D49").Value = "'4.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.45%  '

# Row 50
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.362.04'
$ws.Range("E50").Value = '  +3.15%  '

# Row 51
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").Value = "'2.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.08%  '
